# paises.xlsx -- refresh COVID-19 country statistics + the "updated at" timestamp.
# Source data is sorted by total cases (column B) descending; a handful of
# countries leap-frogged their neighbours with this refresh, so besides updating
# B:H we also rewrite column A for the rows whose country changed position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados a ..." banner (A1): 14:03 -> 15:20 ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 15:20"

# --- Countries that swapped rank (column A only) ---
$ws.Range("A35").Value  = "Paises Bajos"         # was Belgica
$ws.Range("A36").Value  = "Belgica"               # was Paises Bajos
$ws.Range("A78").Value  = "Dinamarca"             # was Australia
$ws.Range("A79").Value  = "Australia"             # was Bosnia y Herzegovina
$ws.Range("A80").Value  = "Bosnia y Herzegovina"  # was Dinamarca
$ws.Range("A207").Value = "Timor Oriental"        # was Santa Lucia
$ws.Range("A208").Value = "Santa Lucia"           # was Timor Oriental

# --- Updated statistics, columns B (Casos totales) .. H (Muertes) ---
# Row 5 (India)
$ws.Range("B5").Value = 6079350
$ws.Range("C5").Value = 6002
$ws.Range("D5").Value = 5017534
$ws.Range("E5").Value = 966210
$ws.Range("G5").Value = 32
$ws.Range("H5").Value = 95606

# Row 19 (Irak)
$ws.Range("B19").Value = 353566
$ws.Range("C19").Value = 4116
$ws.Range("D19").Value = 284784
$ws.Range("E19").Value = 59730
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = 9052

# Row 20 (Arabia Saudita)
$ws.Range("B20").Value = 333648
$ws.Range("C20").Value = 455
$ws.Range("D20").Value = 317846
$ws.Range("E20").Value = 11090
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 4712

# Row 25 (Alemania)
$ws.Range("B25").Value = 286893
$ws.Range("C25").Value = 555
$ws.Range("E25").Value = 26558
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 9535

# Row 32 (Catar)
$ws.Range("B32").Value = 125311
$ws.Range("C32").Value = 227
$ws.Range("D32").Value = 122209
$ws.Range("E32").Value = 2888

# Row 35 (Belgica)
$ws.Range("B35").Value = 114540
$ws.Range("C35").Value = 2914
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("H35").Value = 6380

# Row 36 (Paises Bajos)
$ws.Range("B36").Value = 114179
$ws.Range("C36").Value = 1376
$ws.Range("D36").Value = 19275
$ws.Range("E36").Value = 84924
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 9980

# Row 49 (Bielorrusia)
$ws.Range("B49").Value = 77946
$ws.Range("C49").Value = 337
$ws.Range("D49").Value = 74167
$ws.Range("E49").Value = 2957
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 822

# Row 56 (Barein)
$ws.Range("E56").Value = 6231
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 243

# Row 60 (Uzbekistan)
$ws.Range("B60").Value = 55776
$ws.Range("C60").Value = 456
$ws.Range("D60").Value = 52324
$ws.Range("E60").Value = 2992

# Row 68 (Azerbaiyan)
$ws.Range("B68").Value = 40061
$ws.Range("C68").Value = 38
$ws.Range("D68").Value = 37725
$ws.Range("E68").Value = 1748
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 588

# Row 70 (Estado de Palestina)
$ws.Range("B70").Value = 39121
$ws.Range("C70").Value = 418
$ws.Range("D70").Value = 30220
$ws.Range("E70").Value = 8602
$ws.Range("G70").Value = 8
$ws.Range("H70").Value = 299

# Row 76 (Libia)
$ws.Range("B76").Value = 33213
$ws.Range("C76").Value = 849
$ws.Range("D76").Value = 18518
$ws.Range("E76").Value = 14168
$ws.Range("G76").Value = 7
$ws.Range("H76").Value = 527

# Row 78 (Australia)
$ws.Range("B78").Value = 27072
$ws.Range("C78").Value = 435
$ws.Range("D78").Value = 19942
$ws.Range("E78").Value = 6481
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 649

# Row 79 (Bosnia y Herzegovina)
$ws.Range("B79").Value = 27044
$ws.Range("C79").Value = 4
$ws.Range("D79").Value = 24676
$ws.Range("E79").Value = 1493
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 875

# Row 80 (Dinamarca)
$ws.Range("B80").Value = 26920
$ws.Range("D80").Value = 19746
$ws.Range("E80").Value = 6352
$ws.Range("H80").Value = 822

# Row 92 (Zambia)
$ws.Range("B92").Value = 14660
$ws.Range("C92").Value = 19
$ws.Range("D92").Value = 13821
$ws.Range("E92").Value = 507

# Row 104 (Finlandia)
$ws.Range("E104").Value = 1548
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 345

# Row 105 (Tayikistan)
$ws.Range("B105").Value = 9685
$ws.Range("C105").Value = 39
$ws.Range("D105").Value = 8483
$ws.Range("E105").Value = 1127

# Row 118 (Georgia)
$ws.Range("E118").Value = 3466
$ws.Range("G118").Value = 4
$ws.Range("H118").Value = 32

# Row 149 (Islandia)
$ws.Range("B149").Value = 2663
$ws.Range("C149").Value = 40
$ws.Range("D149").Value = 2161
$ws.Range("E149").Value = 492
